$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their text representation (values are stored as text,
# e.g. "278.96" or "6.66%", not as numbers/percentages) by forcing a Text
# number format before assigning the new value.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "278.96"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "6.66%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "27.49"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.00%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.839"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2.85%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.06346"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "2.48%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "6.944"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "2.97%"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.401"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "7.16%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.8776"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "3.21%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9491"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "4.17%"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "4.76%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.05146"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "1.08%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07329"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03138"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "1.07%"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.09067"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "0.24%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001564"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "1.74%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0006265"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.81%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.006021"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "0.21%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.454"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "0.23%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "0.00%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "3.850"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "-6.57%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04316"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "1.58%"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "0.03%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004300"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "6.02%"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001200"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "-0.06%"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "3.03%"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "3.18%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.006689"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "61.24%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1165"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "4.64%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002157"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "0.67%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01307"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "-1.16%"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "0.97%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "-0.05%"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "857.49%"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "-33.89%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002101"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "-0.05%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002001"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "-0.05%"
